$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$wsVentas = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Row 34 (LINDAO ZUÑIGA BRYAN JOSE / JUAREZ FLORES JORGE WILLIAMS)
$wsVentas.Range("D34").Value = 915.84
$wsVentas.Range("E34").Value = 124.78
$wsVentas.Range("H34").Value = 137.46
$wsVentas.Range("I34").Value = 49.28

# Row 58 (totals "x de 56")
$wsVentas.Range("D58").Value = "1 de 56"
$wsVentas.Range("E58").Value = "2 de 56"
$wsVentas.Range("H58").Value = "1 de 56"
$wsVentas.Range("I58").Value = "2 de 56"

# --- Sheet: VENTA MENSUAL ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# Row 34 (LINDAO ZUÑIGA BRYAN JOSE / JUAREZ FLORES JORGE WILLIAMS)
$wsMensual.Range("F34").Value = 1227.36

# Row 58 (totals)
$wsMensual.Range("F58").Value = 12218.94

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 3 (240X80 PORCELANATO)
$wsCumplimiento.Range("D3").Value = 915.84
$wsCumplimiento.Range("E3").Value = 4076.3432
$wsCumplimiento.Range("F3").Value = 0.1834548059053602

# Row 4 (FREGADEROS DE COCINA)
$wsCumplimiento.Range("D4").Value = 400.37
$wsCumplimiento.Range("E4").Value = -257.867904974973
$wsCumplimiento.Range("F4").Value = 2.809572728946089

# Row 7 (INODOROS)
$wsCumplimiento.Range("D7").Value = 137.46
$wsCumplimiento.Range("E7").Value = 2162.54
$wsCumplimiento.Range("F7").Value = 0.05976521739130435

# Row 8 (LAVABOS)
$wsCumplimiento.Range("D8").Value = 828.6799999999999
$wsCumplimiento.Range("E8").Value = -78.67999999999995
$wsCumplimiento.Range("F8").Value = 1.104906666666667

# Row 19 (TOTAL)
$wsCumplimiento.Range("D19").Value = 12218.94
$wsCumplimiento.Range("E19").Value = 43190.76560036207
$wsCumplimiento.Range("F19").Value = 0.22051985058589
